$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Date with dashes, D, E, F, G, H)
$rows = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 5;  Date = "04-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 6;  Date = "08-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 7;  Date = "11-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 8;  Date = "15-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 9;  Date = "18-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 10; Date = "22-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 11; Date = "25-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 12; Date = "29-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 13; Date = "01-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 14; Date = "05-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 15; Date = "08-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 16; Date = "12-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 17; Date = "15-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 18; Date = "19-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 19; Date = "22-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 20; Date = "26-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 21; Date = "29-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $cellA = $ws.Range("A$n")

    # Force the cell to Text format before assigning the date-like string so
    # Excel's autoconvert doesn't turn e.g. "01-08-2022" into a real date
    # (day <= 12 is ambiguous as dd-mm vs mm-dd). Then restore the default
    # "Normal" style so no stray number-format is left on the cell, matching
    # the original unformatted inline-string cell.
    $cellA.NumberFormat = "@"
    $cellA.Value = $r.Date
    $cellA.Style = "Normal"

    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
}
